$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format first so values like "43.95" or
# "0.691" are written verbatim as strings instead of being coerced into
# floating point numbers (which would introduce binary rounding noise).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '35.531.68'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.893.91'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("D5").Value = '247.28'
$ws.Range("E5").Value = '  -3.56%  '
$ws.Range("D6").Value = '0.691'
$ws.Range("E6").Value = '  -5.67%  '
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").Value = '43.95'
$ws.Range("E8").Value = '  +7.94%  '
$ws.Range("D9").Value = '0.352'
$ws.Range("E9").Value = '  -4.56%  '
$ws.Range("D10").Value = '0.0741'
$ws.Range("E11").Value = '  -1.92%  '
$ws.Range("D12").Value = '13.13'
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").Value = '2.169.44'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").Value = '0.735'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").Value = '4.97'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").Value = '1.899.76'
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").Value = '35.578.79'
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").Value = '73.80'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '0.0₃0825'
$ws.Range("E19").Value = '  -2.63%  '
$ws.Range("D20").Value = '246.84'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").Value = '12.89'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").Value = '4.96'
$ws.Range("E22").Value = '  -2.99%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  +5.00%  '
$ws.Range("D25").Value = '2.18'
$ws.Range("E25").Value = '  -10.26%  '
$ws.Range("D26").Value = '166.48'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '8.52'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").Value = '18.39'
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("D29").Value = '0.127'
$ws.Range("E29").Value = '  -4.40%  '
$ws.Range("D30").Value = '4.128.42'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '1.75'
$ws.Range("E31").Value = '  +6.38%  '
$ws.Range("D32").Value = '4.25'
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").Value = '0.0582'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").Value = '4.22'
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("D36").Value = '0.851'
$ws.Range("E36").Value = '  -6.52%  '
$ws.Range("D37").Value = '2.01'
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("D38").Value = '1.58'
$ws.Range("E38").Value = '  -20.50%  '
$ws.Range("D39").Value = '0.0697'
$ws.Range("E39").Value = '  +7.25%  '
$ws.Range("D40").Value = '17.21'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").Value = '97.82'
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").Value = '0.0216'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").Value = '1.09'
$ws.Range("E43").Value = '  -3.00%  '
$ws.Range("E44").Value = '  -2.89%  '
$ws.Range("D45").Value = '1.293.93'
$ws.Range("E45").Value = '  -3.00%  '
$ws.Range("D46").Value = '0.0817'
$ws.Range("E46").Value = '  +8.25%  '
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '2.74'
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("D49").Value = '12.06'
$ws.Range("E49").Value = '  +2.97%  '
$ws.Range("D50").Value = '43.33'
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("D51").Value = '6.33'
$ws.Range("E51").Value = '  -5.82%  '

# Restore the default (unstyled) cell style on the Price column so the
# saved XML doesn't carry a leftover text-format style attribute that
# wasn't present in the original workbook.
$ws.Range("D2:D51").Style = "Normal"
